$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the duplicated "Play Black Widow Slot Game for Free |
#    Review" paragraph that used to sit right before the final
#    (italic) paragraph, and turn that final paragraph's text into
#    the new image-generation prompt.
# ------------------------------------------------------------------
$headingText = "Play Black Widow Slot Game for Free | Review"

$dupIndex = -1
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $ptxt = $d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($ptxt -eq $headingText) {
        $dupIndex = $i
        break
    }
}
if ($dupIndex -ge 2) {
    $d.Paragraphs($dupIndex).Range.Delete()
}

$oldBlurb = "Read our review of Black Widow slot game, available to play for free. Get a chance to win huge prizes with the game's free spins bonus round."
$newPrompt = "Create a feature image for Black Widow slot game featuring a happy Maya warrior with glasses in a cartoon style. The image should feature the warrior holding a spider and standing in front of a spider web. The background should be dark with cobwebs on the corners to match the theme of the game. The warrior should be dressed in a black jumpsuit with a red hourglass symbol on the chest and his/her arms folded in front. The image should be eye-catching with vibrant colors to attract players to the game."

$d.Content.Find.Execute($oldBlurb, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2) | Out-Null

# ------------------------------------------------------------------
# 2) Insert a new "Meta description" paragraph right after the
#    document's opening Heading1 title paragraph.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"

$boldLabel = "Meta description"
$restOfLine = ": Read our review of Black Widow slot game, available to play for free. Get a chance to win huge prizes with the game's free spins bonus round."

$metaStart = $metaPara.Range.Start
$insertionRange = $d.Range($metaStart, $metaStart)
$insertionRange.InsertAfter($boldLabel + $restOfLine)

$boldRange = $d.Range($metaStart, $metaStart + $boldLabel.Length)
$boldRange.Font.Bold = 1
